$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the two retired test-scenario rows -------------------------------
# "RSVP - Basic" / "Member cancels RSVP (back to no response)" (row 29)
# "Events" / "Host assigns another host" (row 25)
# Delete bottom-most row first so the other row number stays valid.
$ws.Rows.Item(29).Delete()
$ws.Rows.Item(25).Delete()

# --- Mark Desktop/Mobile Web/Flutter coverage ("X") for the Events rows that --
# --- now have it, and for the two remaining RSVP - Basic rows ----------------
$eventsRows = @(20, 21, 23, 24, 26, 27)
foreach ($r in $eventsRows) {
    $ws.Range("C$r").Value = "X"
    $ws.Range("D$r").Value = "X"
    $ws.Range("E$r").Value = "X"
}

# --- Restore the selection/active cell used in the saved view ----------------
$ws.Range("A33").Select()
